$wb = $excel.ActiveWorkbook

# ---- Metadata sheet updates ----
$meta = $wb.Worksheets.Item("Metadata")

# Date property (row 8)
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# FHIR Version property (row 15)
$meta.Range("B15").Value = "4.0.1"

# ---- Elements sheet updates ----
$els = $wb.Worksheets.Item("Elements")

# Row 2 ("Extension"): update ele-1 constraint text (drop the "unless an empty Parameters..." clause)
$els.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 ("Extension.id"): Type(s) changes from "id" to "string"
$els.Range("K3").Value = "string" + [char]10

# Row 6 ("Extension.value[x]"): Definition text changes R4B -> R4 in the extensibility link
$els.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."

$wb.Save()
